$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; everything shifts one column right.
$ws.Columns("A").Insert()

# New column A header (row 3) — "Match ID", bold (no border) — matches the
# style used for the other un-bordered bold header cells.
$ws.Range("A3").Value = "Match ID"
$ws.Range("A3").Font.Bold = $true

# Data rows 4-19: Match ID = 26, same bold style as the header.
$ws.Range("A4:A19").Value = 26
$ws.Range("A4:A19").Font.Bold = $true

# Row 20 (the hidden totals row) keeps the default (unstyled) cell format.
$ws.Range("A20").Value = 26

# Restore the selection to the new Match ID column's data cells.
$ws.Range("A3:A19").Select() | Out-Null
